# fix embryo import dan test logic
# - remove the "ID Worker" column (shifts all following columns left by one)
# - change "Jlh Botol" values for row 2/3
# - add a trailing "<end>" marker row
# - style header row (bold, fill, border, centered) and center the data cells
# - adjust zoom / selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Remove the "ID Worker" column (column B). Everything to the
#    right (Sub Culture, Jlh Botol, Tgl Botol, ...) shifts left.
# ------------------------------------------------------------------
$ws.Columns.Item(2).Delete()

# ------------------------------------------------------------------
# 2. Update the "Jlh Botol" data values (now column C).
# ------------------------------------------------------------------
$ws.Range("C2").Value = 100
$ws.Range("C3").Value = 200

# ------------------------------------------------------------------
# 3. Append the "<end>" marker row.
# ------------------------------------------------------------------
$ws.Range("A4").Value = "<end>"

# ------------------------------------------------------------------
# 4. Style the header row: bold font, light grey fill, thin border,
#    centered text.
# ------------------------------------------------------------------
$headerRange = $ws.Range("A1:H1")
$headerRange.HorizontalAlignment = -4108
$headerRange.Borders.LineStyle = 1
$headerRange.Interior.Pattern = 1
$headerRange.Interior.ThemeColor = 2
$headerRange.Interior.TintAndShade = -0.15
$headerRange.Font.Bold = $true

# ------------------------------------------------------------------
# 5. Center the numeric data cells (everything except the date column).
# ------------------------------------------------------------------
$ws.Range("A2:C3").HorizontalAlignment = -4108
$ws.Range("E2:H3").HorizontalAlignment = -4108

# ------------------------------------------------------------------
# 6. Center the date cells (column D) while keeping their original
#    (non-default) font.
# ------------------------------------------------------------------
$dateRange = $ws.Range("D2:D3")
$dateRange.HorizontalAlignment = -4108
$dateRange.Font.ThemeFont = 1
$dateRange.Font.Size = 11

# ------------------------------------------------------------------
# 7. Center the "<end>" cell, also keeping a non-default font.
# ------------------------------------------------------------------
$endRange = $ws.Range("A4")
$endRange.HorizontalAlignment = -4108
$endRange.Font.ThemeFont = 1
$endRange.Font.Size = 11

# ------------------------------------------------------------------
# 8. Update view: zoom 160%, selection on D2.
# ------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 160
$ws.Range("D2").Select()
